$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1044271, 348591, 73075, 1465937)
    3 = @(716689, 327971, 11024, 1055684)
    4 = @(1236227, 919531, 52972, 2208730)
    5 = @(1124792, 628490, 154957, 1908239)
    6 = @(1178695, 481616, 16756, 1677067)
    7 = @(296863, 85904, 38284, 421051)
    8 = @(450695, 18116, 4023, 472834)
    9 = @(758072, 413298, 49352, 1220722)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
}
